# Applies the cryptos-list refresh described by the commit diff.
# Each entry: cell reference, new text value, and whether the value must be
# forced to Text (otherwise Excel's COM layer would auto-coerce a clean
# decimal string like "540.88" into a numeric cell, which the source diff
# does not want -- every cell in this sheet is a plain text cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '60.430.60'; ForceText = $false },
    @{ Cell = 'E2'; Value = '  -1.63%  '; ForceText = $false },
    @{ Cell = 'D3'; Value = '2.335.19'; ForceText = $false },
    @{ Cell = 'E3'; Value = '  -4.58%  '; ForceText = $false },
    @{ Cell = 'E4'; Value = '  +0.12%  '; ForceText = $false },
    @{ Cell = 'D5'; Value = '540.88'; ForceText = $true },
    @{ Cell = 'E5'; Value = '  -1.16%  '; ForceText = $false },
    @{ Cell = 'D6'; Value = '136.15'; ForceText = $true },
    @{ Cell = 'E6'; Value = '  -7.05%  '; ForceText = $false },
    @{ Cell = 'E7'; Value = '  +0.09%  '; ForceText = $false },
    @{ Cell = 'D8'; Value = '0.522'; ForceText = $true },
    @{ Cell = 'E8'; Value = '  -10.54%  '; ForceText = $false },
    @{ Cell = 'D9'; Value = '2.335.07'; ForceText = $false },
    @{ Cell = 'E9'; Value = '  -4.43%  '; ForceText = $false },
    @{ Cell = 'E10'; Value = '  -1.93%  '; ForceText = $false },
    @{ Cell = 'D11'; Value = '0.155'; ForceText = $true },
    @{ Cell = 'E11'; Value = '  +0.08%  '; ForceText = $false },
    @{ Cell = 'D12'; Value = '5.27'; ForceText = $true },
    @{ Cell = 'E12'; Value = '  -2.19%  '; ForceText = $false },
    @{ Cell = 'D13'; Value = '0.339'; ForceText = $true },
    @{ Cell = 'E13'; Value = '  -3.02%  '; ForceText = $false },
    @{ Cell = 'D14'; Value = '24.42'; ForceText = $true },
    @{ Cell = 'E14'; Value = '  -6.02%  '; ForceText = $false },
    @{ Cell = 'D15'; Value = '2.756.02'; ForceText = $false },
    @{ Cell = 'E15'; Value = '  -4.41%  '; ForceText = $false },
    @{ Cell = 'D16'; Value = '60.467.49'; ForceText = $false },
    @{ Cell = 'E16'; Value = '  -1.29%  '; ForceText = $false },
    @{ Cell = 'E17'; Value = '  -4.20%  '; ForceText = $false },
    @{ Cell = 'D18'; Value = '2.335.18'; ForceText = $false },
    @{ Cell = 'E18'; Value = '  -4.38%  '; ForceText = $false },
    @{ Cell = 'D19'; Value = '10.49'; ForceText = $true },
    @{ Cell = 'E19'; Value = '  -3.84%  '; ForceText = $false },
    @{ Cell = 'D20'; Value = '316.50'; ForceText = $true },
    @{ Cell = 'E20'; Value = '  -0.37%  '; ForceText = $false },
    @{ Cell = 'D21'; Value = '4.05'; ForceText = $true },
    @{ Cell = 'E21'; Value = '  -2.31%  '; ForceText = $false },
    @{ Cell = 'D22'; Value = '6.54'; ForceText = $true },
    @{ Cell = 'E22'; Value = '  -5.70%  '; ForceText = $false },
    @{ Cell = 'D23'; Value = '0.999'; ForceText = $true },
    @{ Cell = 'E23'; Value = '  -0.09%  '; ForceText = $false },
    @{ Cell = 'D24'; Value = '1.84'; ForceText = $true },
    @{ Cell = 'E24'; Value = '  -2.23%  '; ForceText = $false },
    @{ Cell = 'D25'; Value = '62.89'; ForceText = $true },
    @{ Cell = 'E25'; Value = '  -1.09%  '; ForceText = $false },
    @{ Cell = 'D26'; Value = '8.57'; ForceText = $true },
    @{ Cell = 'E26'; Value = '  +10.67%  '; ForceText = $false },
    @{ Cell = 'D27'; Value = '0.998'; ForceText = $true },
    @{ Cell = 'E27'; Value = '  +0.14%  '; ForceText = $false },
    @{ Cell = 'D28'; Value = '2.449.29'; ForceText = $false },
    @{ Cell = 'E28'; Value = '  -4.42%  '; ForceText = $false },
    @{ Cell = 'B29'; Value = 'PEPE'; ForceText = $false },
    @{ Cell = 'C29'; Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'; ForceText = $false },
    @{ Cell = 'D29'; Value = '0.0₃0888'; ForceText = $false },
    @{ Cell = 'E29'; Value = '  -9.10%  '; ForceText = $false },
    @{ Cell = 'B30'; Value = 'InternetComputer(DFINITY)'; ForceText = $false },
    @{ Cell = 'C30'; Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; ForceText = $false },
    @{ Cell = 'D30'; Value = '7.91'; ForceText = $true },
    @{ Cell = 'E30'; Value = '  -3.84%  '; ForceText = $false },
    @{ Cell = 'E31'; Value = '  -5.75%  '; ForceText = $false },
    @{ Cell = 'D32'; Value = '498.32'; ForceText = $true },
    @{ Cell = 'E32'; Value = '  -7.66%  '; ForceText = $false },
    @{ Cell = 'E33'; Value = '  -1.90%  '; ForceText = $false },
    @{ Cell = 'E34'; Value = '  -5.15%  '; ForceText = $false },
    @{ Cell = 'E35'; Value = '  -3.78%  '; ForceText = $false },
    @{ Cell = 'E36'; Value = '  +0.22%  '; ForceText = $false },
    @{ Cell = 'D37'; Value = '4.59'; ForceText = $true },
    @{ Cell = 'E37'; Value = '  -4.44%  '; ForceText = $false },
    @{ Cell = 'D38'; Value = '0.372'; ForceText = $true },
    @{ Cell = 'E38'; Value = '  -1.76%  '; ForceText = $false },
    @{ Cell = 'D39'; Value = '18.23'; ForceText = $true },
    @{ Cell = 'E39'; Value = '  +0.11%  '; ForceText = $false },
    @{ Cell = 'D40'; Value = '5.21'; ForceText = $true },
    @{ Cell = 'E40'; Value = '  -9.13%  '; ForceText = $false },
    @{ Cell = 'E41'; Value = '  +1.69%  '; ForceText = $false },
    @{ Cell = 'D42'; Value = '0.999'; ForceText = $true },
    @{ Cell = 'E42'; Value = '  -0.07%  '; ForceText = $false },
    @{ Cell = 'D43'; Value = '137.48'; ForceText = $true },
    @{ Cell = 'E43'; Value = '  -1.75%  '; ForceText = $false },
    @{ Cell = 'D44'; Value = '40.09'; ForceText = $true },
    @{ Cell = 'E44'; Value = '  -0.23%  '; ForceText = $false },
    @{ Cell = 'D45'; Value = '140.83'; ForceText = $true },
    @{ Cell = 'E45'; Value = '  -0.59%  '; ForceText = $false },
    @{ Cell = 'E46'; Value = '  -9.04%  '; ForceText = $false },
    @{ Cell = 'D47'; Value = '3.52'; ForceText = $true },
    @{ Cell = 'E47'; Value = '  -1.97%  '; ForceText = $false },
    @{ Cell = 'D48'; Value = '0.0509'; ForceText = $true },
    @{ Cell = 'E48'; Value = '  -4.70%  '; ForceText = $false },
    @{ Cell = 'D49'; Value = '19.33'; ForceText = $true },
    @{ Cell = 'E49'; Value = '  -10.31%  '; ForceText = $false },
    @{ Cell = 'E50'; Value = '  -3.69%  '; ForceText = $false },
    @{ Cell = 'D51'; Value = '0.0896'; ForceText = $true },
    @{ Cell = 'E51'; Value = '  -3.42%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        # Round-trip through a Text number format so the literal digits are
        # stored as a string (matching the inline-string cells in the source
        # file), then drop the temporary format so no stray style lingers on
        # the cell.
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.ClearFormats()
    } else {
        $rng.Value = $u.Value
    }
}
